$wb = $excel.ActiveWorkbook

# Rename sheets with "IA_" prefix
foreach ($ws in $wb.Worksheets) {
    $ws.Name = "IA_" + $ws.Name
}

# Make the "IA_Attainment" sheet the active/selected sheet (was Demographics, now Attainment)
$wb.Worksheets.Item("IA_Attainment").Activate()
